$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 2. Data reporter section updates
$ws.Range("B10").Value = "www.stat.gov.kg"
$ws.Range("B9").Value = "(0312) 62 56 07"
$ws.Range("B7").Value = "Mambetaliev T.A."
$ws.Range("B6").Value = "National Statistical Committee of the Kyrgyz Republic (Department of Digital Development and Sustainable Development Statistics)"

# The updated contact-person / phone cells were retyped using a Cyrillic
# keyboard layout, which made Excel tag their runs with a second Calibri
# font record (charset 204). Re-assert the font on B7/B9 (no-wrap cells)
# and then B6 (wrap cell) to reproduce that extra font + style split.
$ws.Range("B7").Font.Name = "Calibri"
$ws.Range("B9").Font.Name = "Calibri"
$ws.Range("B6").Font.Name = "Calibri"

# Move the active selection to B8
$ws.Range("B8").Select()
